# Hortaliza, Feria Lagunitas de Puerto Montt - Ajo
# Insert a new weekly price record at row 378 (pushing the existing
# rows 378-393 down to 379-394), matching the author's weekly update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 378, shifting rows 378:393 down to 379:394.
$ws.Rows.Item(378).Insert()

# Populate the new row with the latest week's data.
$ws.Range("A378").Value = 4
$ws.Range("B378").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C378").Value = "Los Lagos"
$ws.Range("D378").Value = 44939
$ws.Range("E378").Value = 10
$ws.Range("F378").Value = 100112003
$ws.Range("G378").Value = "Ajo"
$ws.Range("H378").Value = "Chino"
$ws.Range("I378").Value = "Primera"
$ws.Range("J378").Value = 240
$ws.Range("K378").Value = 18000
$ws.Range("L378").Value = 19000
$ws.Range("M378").Value = 18500
$ws.Range("N378").Value = "$/caja 10 kilos"
$ws.Range("O378").Value = "China"
$ws.Range("P378").Value = 1850
$ws.Range("Q378").Value = 10
$ws.Range("R378").Value = "Hortaliza"
